$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 39

# Columns A, B, D, E, F hold non-numeric-looking text, so plain assignment
# keeps them as text cells (matching the inlineStr cells elsewhere in the
# sheet).
$ws.Cells.Item($row, 1).Value = "2024-09-25T18:06:40Z"
$ws.Cells.Item($row, 2).Value = "temperature"
$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "N/A"
$ws.Cells.Item($row, 6).Value = "N/A"

# Column C ("25") looks numeric, so Excel would normally coerce it to a
# number. Temporarily force a text format while assigning it, then restore
# the cell's style to Normal so no stray number-format style lingers on the
# cell (matching the rest of the sheet, which has no explicit style index).
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "25"
$ws.Cells.Item($row, 3).Style = "Normal"
